$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.442.42"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "1.823.36"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.46"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4586"
$ws.Range("E7").Value = "  -1.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3807"
$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.35"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07888"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9677"
$ws.Range("E11").Value = "  -3.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.01"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.842.34"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.878"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.040"
$ws.Range("E15").Value = "  -2.48%  "

$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.73"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06631"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001025"
$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").Value = "27.442.29"
$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.334"
$ws.Range("E23").Value = "  -1.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.79"
$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "2.046.38"
$ws.Range("E26").Value = "  -1.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.65"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.33"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.055"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.261"
$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.17"
$ws.Range("E31").Value = "  -2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9430"
$ws.Range("E32").Value = "  -2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09318"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.593"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.234"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.320"
$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05918"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02177"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.159"
$ws.Range("E39").Value = "  -2.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.021"
$ws.Range("E40").Value = "  -2.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5760"
$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1828"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.260"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.00"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5435"
$ws.Range("E46").Value = "  -3.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.865"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06604"
$ws.Range("E48").Value = "  -2.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.29"
$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.039"
$ws.Range("E51").Value = "  -1.44%  "
